# Replace the text-encoded measurement values with plain numbers, add the
# missing RA300hr/60oC readings for rows 6-7, and move the active selection
# the way the author left it (I11 instead of I12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Voltage (Vp-p), Frequency (Hz) and Time (hr) columns: drop the unit
# suffixes that were stored as text ("14Vp-p", "60Hz ", "72h", "RA300hr")
# and store the bare numeric readings instead - for every data row.
$ws.Range("F2:F7").Value = 14
$ws.Range("G2:G7").Value = 60

$ws.Range("H2:H5").Value = 72
$ws.Range("H6:H7").Value = 300

# Temperature column: rows 6-7 used to read the text "60oC"; now a bare 60.
$ws.Range("I6:I7").Value = 60

# Vender / File source were blank for the 300hr/60oC rows; fill with T.B.D.
$ws.Range("J6:J7").Value = "T.B.D"
$ws.Range("K6:K7").Value = "T.B.D"

# Author's last selection before saving.
[void]$ws.Range("I11").Select()
